$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A ("TabName") before the existing data, shifting
# everything one column to the right.
$ws.Columns.Item(1).Insert()

# Insert a new row 3 (for the "FilesTab" entry) after the existing row 2.
$ws.Rows.Item(3).Insert()

$casesTabQuery = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 WHERE a.arm_id IN ['Q']
OPTIONAL MATCH (f:file)-[*]->(c)
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@

$statQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
WHERE a.arm_id IN ['Q']
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

$filesTabQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
WHERE a.arm_id IN ['Q']
WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@

# Assign the brand-new shared strings in the same order the authoring
# session originally created them in (TabName, CasesTab, FilesTab,
# FilesTab-query, CasesTab-query, StatQuery) so the sharedStrings table
# comes out in the same order as the saved workbook.
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"
$ws.Range("A3").Value = "FilesTab"
$ws.Range("B3").Value = $filesTabQuery
$ws.Range("B2").Value = $casesTabQuery
$ws.Range("C2").Value = $statQuery

# ---- Remaining, already-existing values ----
$ws.Range("B1").Value = "query"
$ws.Range("C1").Value = "StatQuery"
$ws.Range("D1").Value = "dbExcel"
$ws.Range("E1").Value = "WebExcel"

$ws.Range("D2").Value = "TC01_Trials_Filter_TrialArm-Q_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC01_Trials_Filter_TrialArm-Q_WebData.xlsx"

$ws.Range("C3").Value = $statQuery
$ws.Range("D3").Value = "TC01_Trials_Filter_TrialArm-Q_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC01_Trials_Filter_TrialArm-Q_WebData.xlsx"

# ---- Formatting: wrap text on the query cells (reuses existing wrap style) ----
$ws.Range("B2").WrapText = $true
$ws.Range("C2").WrapText = $true
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# ---- Row heights ----
$ws.Rows.Item(2).RowHeight = 195
$ws.Rows.Item(3).RowHeight = 409.5

# ---- Column widths ----
$ws.Columns.Item(1).ColumnWidth = 8.85546875
$ws.Columns.Item(2).ColumnWidth = 75.85546875
$ws.Columns.Item(3).ColumnWidth = 75.85546875
$ws.Columns.Item(4).ColumnWidth = 70.28515625
$ws.Columns.Item(5).ColumnWidth = 28.5703125

# ---- View / selection ----
$ws.Range("B3").Select()
